# Apply latest cryptos price/volume(1h) snapshot
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.747.37"
$ws.Range("E2").Value = "  -0.22%  "

$ws.Range("D3").Value = "1.632.40"
$ws.Range("E3").Value = "  -0.22%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.05%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -0.85%  "

$ws.Range("E9").Value = "  -1.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0792"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.97%  "

$ws.Range("E12").Value = "  +0.45%  "

$ws.Range("D13").Value = "1.858.08"
$ws.Range("E13").Value = "  -0.19%  "

$ws.Range("D14").Value = "1.625.92"
$ws.Range("E14").Value = "  -0.65%  "

$ws.Range("E15").Value = "  +0.23%  "

$ws.Range("D16").Value = "0.0₃0765"
$ws.Range("E16").Value = "  -1.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.05%  "

$ws.Range("D18").Value = "25.772.77"
$ws.Range("E18").Value = "  -0.15%  "

$ws.Range("E19").Value = "  -0.06%  "

$ws.Range("E20").Value = "  +0.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "192.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.84%  "

$ws.Range("E23").Value = "  +1.77%  "

$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("E25").Value = "  +2.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "143.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.82%  "

$ws.Range("E27").Value = "  +1.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.56%  "

$ws.Range("E29").Value = "  -0.34%  "

$ws.Range("E30").Value = "  -0.05%  "

$ws.Range("E31").Value = "  -0.97%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.18%  "

$ws.Range("E33").Value = "  -0.91%  "

$ws.Range("E35").Value = "  -0.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.905"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.82%  "

$ws.Range("D37").Value = "1.132.42"
$ws.Range("E37").Value = "  +2.15%  "

$ws.Range("E38").Value = "  -1.92%  "

$ws.Range("E39").Value = "  -1.26%  "

$ws.Range("E40").Value = "  -1.06%  "

$ws.Range("E41").Value = "  +0.12%  "

$ws.Range("E42").Value = "  +0.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.66%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.799"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.31%  "

$ws.Range("D46").Value = "1.767.08"
$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("E47").Value = "  +0.55%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.30%  "

$ws.Range("E50").Value = "  -0.21%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.24%  "
